$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row with the new variable entry
$ws.Range("A9").Value = 1746359037
$ws.Range("B9").Value = "update"
$ws.Range("C9").Value = "variable"
$ws.Range("D9").Value = "dep_sante___variable_3"
$ws.Range("F9").Value = "type"
$ws.Range("G9").Value = "integer"
$ws.Range("H9").Value = "string"
